$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three same-position value corrections in the existing August block (rows 2-13 keep
# their row numbers; only total_venda for days 4, 12 and 13 changed).
$ws.Range("B3").Value  = 30365.46   # Dia 4  (Ago/2025)
$ws.Range("B9").Value  = 115305.79  # Dia 12 (Ago/2025)
$ws.Range("B10").Value = 26032      # Dia 13 (Ago/2025)

# A new record (Dia 19, Ago/2025) was inserted after the existing row 13, pushing every
# following record down by one row (old row 14 -> new row 15, ..., old row 78 -> new row 79).
# Row-insert operations are avoided here (they perturb unrelated numeric literals), so the
# whole block from row 14 to row 79 is rewritten in place with its final values instead.
$rows = @(
    @(19, 10467.2, 8, 2025, "08/2025"),
    @(1, 18112.93, 7, 2025, "07/2025"),
    @(2, 27735.81, 7, 2025, "07/2025"),
    @(3, 16620.34, 7, 2025, "07/2025"),
    @(4, 27892.39, 7, 2025, "07/2025"),
    @(7, 24182.46, 7, 2025, "07/2025"),
    @(8, 27515.9, 7, 2025, "07/2025"),
    @(9, 15367.22, 7, 2025, "07/2025"),
    @(10, 21120.89, 7, 2025, "07/2025"),
    @(11, 14249.1, 7, 2025, "07/2025"),
    @(14, 18544.3, 7, 2025, "07/2025"),
    @(15, 48059.45, 7, 2025, "07/2025"),
    @(16, 15067.65, 7, 2025, "07/2025"),
    @(17, 12214.2, 7, 2025, "07/2025"),
    @(18, 7505.85, 7, 2025, "07/2025"),
    @(21, 499270.27, 7, 2025, "07/2025"),
    @(22, 9220.86, 7, 2025, "07/2025"),
    @(23, 12454.83, 7, 2025, "07/2025"),
    @(24, 5000.2, 7, 2025, "07/2025"),
    @(25, 7210.32, 7, 2025, "07/2025"),
    @(28, 20062.57, 7, 2025, "07/2025"),
    @(29, 25018.04, 7, 2025, "07/2025"),
    @(30, 59573.83, 7, 2025, "07/2025"),
    @(31, 27720.5, 7, 2025, "07/2025"),
    @(2, 45067.52, 6, 2025, "06/2025"),
    @(3, 35136.72, 6, 2025, "06/2025"),
    @(4, 11872.2, 6, 2025, "06/2025"),
    @(5, 12749.74, 6, 2025, "06/2025"),
    @(6, 12978.82, 6, 2025, "06/2025"),
    @(9, 11059.45, 6, 2025, "06/2025"),
    @(10, 4973.1, 6, 2025, "06/2025"),
    @(11, 19810.51, 6, 2025, "06/2025"),
    @(12, 18059.33, 6, 2025, "06/2025"),
    @(13, 9230.38, 6, 2025, "06/2025"),
    @(16, 24913.77, 6, 2025, "06/2025"),
    @(17, 14906.35, 6, 2025, "06/2025"),
    @(18, 14248.93, 6, 2025, "06/2025"),
    @(20, 5940.05, 6, 2025, "06/2025"),
    @(21, 119.13, 6, 2025, "06/2025"),
    @(23, 46214.09, 6, 2025, "06/2025"),
    @(24, 15253.9, 6, 2025, "06/2025"),
    @(25, 14231.97, 6, 2025, "06/2025"),
    @(26, 24113.52, 6, 2025, "06/2025"),
    @(27, 9530.56, 6, 2025, "06/2025"),
    @(30, 114294.26, 6, 2025, "06/2025"),
    @(2, 20463.44, 5, 2025, "05/2025"),
    @(5, 29720.49, 5, 2025, "05/2025"),
    @(6, 19666.55, 5, 2025, "05/2025"),
    @(7, 21883.41, 5, 2025, "05/2025"),
    @(8, 45418.95, 5, 2025, "05/2025"),
    @(9, 39575.07, 5, 2025, "05/2025"),
    @(12, 15499.7, 5, 2025, "05/2025"),
    @(13, 10155.95, 5, 2025, "05/2025"),
    @(14, 34413.69, 5, 2025, "05/2025"),
    @(15, 28325.85, 5, 2025, "05/2025"),
    @(16, 12000.74, 5, 2025, "05/2025"),
    @(19, 17756.75, 5, 2025, "05/2025"),
    @(20, 23593.88, 5, 2025, "05/2025"),
    @(21, 9475.47, 5, 2025, "05/2025"),
    @(22, 27766.95, 5, 2025, "05/2025"),
    @(23, 11639.15, 5, 2025, "05/2025"),
    @(26, 28728.43, 5, 2025, "05/2025"),
    @(27, 15807.63, 5, 2025, "05/2025"),
    @(28, 24406.06, 5, 2025, "05/2025"),
    @(29, 30466.72, 5, 2025, "05/2025"),
    @(30, 23720.49, 5, 2025, "05/2025")
)

$r = 14
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
